$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new driver entry (was 23.90.0.2, now 23.10.0.8) with updated counts
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 98.8

# Row 4: driver shifted (was 23.60.1.2, now 23.90.0.2) with updated counts
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2"
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 374

# Row 5: totals updated
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 380

# Row 13: clear the Driver Vintage date value
$ws.Range("E13").Value = ""

# Row 15: updated sample count
$ws.Range("B15").Value = 265400
